# "Reworking the admin part"
# - Drop the redundant "Группа " ("Group ") prefix from both sheet tab names.
# - Switch the active/selected tab from the first sheet ("А-1-25") to the
#   second sheet ("Б-1-25").

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "А-1-25"
$ws2.Name = "Б-1-25"

# Make the second sheet the active one (moves tabSelected / activeTab).
$ws2.Activate()
